# Commit: "remove column from alcohol data"
# Sheet1 had an extra numeric column (M) that duplicated/obsoleted the data
# in the following column (N). Delete column M so the old column N shifts
# left and becomes the new column M (dimension goes from A1:N119 to A1:M119).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").EntireColumn.Delete()
